# Added arabic (and french) lang rows to the identity_schema master-data sheet.
#
# Row 2 (the "eng" row) is used as the template: every new row duplicates
# the formatting/content of row 2 except for the lang_code (column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Copy-CellFormat($srcRow, $srcCol, $dstRow, $dstCol) {
    $ws.Cells.Item($srcRow, $srcCol).Copy() | Out-Null
    $ws.Cells.Item($dstRow, $dstCol).PasteSpecial($xlPasteFormats) | Out-Null
}

function Add-LangRow($rowNum, $langCode) {
    # Column A: lang_code (new shared string, e.g. "fra" / "ara")
    $ws.Cells.Item($rowNum, 1).Value = $langCode
    Copy-CellFormat 2 1 $rowNum 1

    # Column B: id -> numeric 1001 (column is text-formatted, so force
    # General before/after the write to keep the stored type numeric)
    $ws.Cells.Item($rowNum, 2).NumberFormat = "General"
    $ws.Cells.Item($rowNum, 2).Value = 1001
    $ws.Cells.Item($rowNum, 2).NumberFormat = "@"

    # Column C: id_version -> numeric 0.1
    $ws.Cells.Item($rowNum, 3).Value = 0.1
    Copy-CellFormat 2 3 $rowNum 3

    # Column D: title
    $ws.Cells.Item($rowNum, 4).Value = $ws.Cells.Item(2, 4).Value()
    Copy-CellFormat 2 4 $rowNum 4

    # Column E: description
    $ws.Cells.Item($rowNum, 5).Value = $ws.Cells.Item(2, 5).Value()
    Copy-CellFormat 2 5 $rowNum 5

    # Column F: schema_json (large text blob)
    $ws.Cells.Item($rowNum, 6).Value = $ws.Cells.Item(2, 6).Value()
    Copy-CellFormat 2 6 $rowNum 6

    # Column G: status_code ("PUBLISHED")
    $ws.Cells.Item($rowNum, 7).Value = $ws.Cells.Item(2, 7).Value()
    Copy-CellFormat 2 7 $rowNum 7

    # Column H: add_props ("FALSE" as text, not boolean) - use the quote
    # prefix trick to force text, then restore the clean number format by
    # pasting row 2's format over it.
    $ws.Cells.Item($rowNum, 8).Value = "'FALSE"
    Copy-CellFormat 2 8 $rowNum 8

    # Column I: effective_from ("now()")
    $ws.Cells.Item($rowNum, 9).Value = $ws.Cells.Item(2, 9).Value()
    Copy-CellFormat 2 9 $rowNum 9

    # Column J: is_active ("TRUE" as text, not boolean)
    $ws.Cells.Item($rowNum, 10).Value = "'TRUE"
    Copy-CellFormat 2 10 $rowNum 10

    # Row height matches the other data rows (auto height driven by the
    # wrapped schema_json column in the real workbook).
    $ws.Rows.Item($rowNum).RowHeight = 1120.85
}

Add-LangRow 3 "fra"
Add-LangRow 4 "ara"

# Move the selection / active cell to A4, like the source workbook.
$ws.Range("A4").Select() | Out-Null
